$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("B-T-x")
try {
  $excel.ActivePrinter = "Microsoft Print to PDF"
} catch {
  Write-Output "ERR1: $_"
}
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
